$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rental_item table: new index "idx_condition" on (pickup_condition, return_condition) ---
$ws.Range("F102").Value = "pickup_condition, return_condition"
$ws.Range("A102").Value = "idx_condition"
$ws.Range("D102").Value = "x"

# --- employee table: new index "idx_role_client" on (role, id_client) ---
$ws.Range("F31").Value = "role, id_client"
$ws.Range("A31").Value = "idx_role_client"
$ws.Range("D31").Value = "x"

# --- rental table: new index "idx_status_start" on (status, start_date) ---
$ws.Range("F83").Value = "status, start_date"
$ws.Range("A83").Value = "idx_status_start"
$ws.Range("D83").Value = "x"

# --- rental_item table: new index "idx_rentalitem_equipment" on (id_rental, id_equipment) ---
$ws.Range("F103").Value = "id_rental, id_equipment"
$ws.Range("A103").Value = "idx_rentalitem_equipment"
$ws.Range("D103").Value = "x"

# --- view state: scroll back to top and move selection to H6 ---
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("H6").Select()
